$wb = $excel.ActiveWorkbook

# Sheet 1: ALC
$ws = $wb.Worksheets.Item(1)
$ws.Range("H123").Value = 25976.2
$ws.Range("H133").Value = 25800
$ws.Range("H136").Value = 0
$ws.Range("H58").Value = 1287.8
$ws.Range("I58").Value = 475.8
$ws.Range("J123").Value = 25976.2
$ws.Range("J133").Value = 25800
$ws.Range("J136").Value = 0
$ws.Range("J58").Value = 2099.8
$ws.Range("K58").Value = 1427.4
$ws.Range("L123").Value = 25976.2
$ws.Range("L133").Value = 25800
$ws.Range("L136").Value = 0
$ws.Range("L58").Value = 6299.400000000001
$ws.Range("M58").Value = -1277.4
$ws.Range("N123").Value = -35776.2
$ws.Range("N133").Value = -35920
$ws.Range("N136").ClearContents()
$ws.Range("N58").Value = -6599.400000000001

# Sheet 2: ARM
$ws = $wb.Worksheets.Item(2)
$ws.Range("H100").Value = 45177.5
$ws.Range("H122").Value = 2565644.8
$ws.Range("H123").Value = 43857
$ws.Range("H131").Value = 46000
$ws.Range("H24").Value = 45177.5
$ws.Range("H32").Value = 8841.912
$ws.Range("H45").Value = 8680.733
$ws.Range("H63").Value = 100002696
$ws.Range("H66").Value = 100002696
$ws.Range("I122").Value = 2850494
$ws.Range("I32").Value = 8767.237999999999
$ws.Range("I45").Value = 11110.091
$ws.Range("I63").Value = 111113720
$ws.Range("I66").Value = 111113720
$ws.Range("J100").Value = 45177.5
$ws.Range("J122").Value = 2000
$ws.Range("J123").Value = 43857
$ws.Range("J131").Value = 46000
$ws.Range("J24").Value = 45177.5
$ws.Range("J32").Value = 8962.538
$ws.Range("J63").Value = 3500
$ws.Range("J66").Value = 3500
$ws.Range("K122").Value = 8551482
$ws.Range("K32").Value = 8767.237999999999
$ws.Range("K45").Value = 11110.091
$ws.Range("K63").Value = 111113720
$ws.Range("K66").Value = 555568600
$ws.Range("L100").Value = 45177.5
$ws.Range("L122").Value = 6000
$ws.Range("L123").Value = 43857
$ws.Range("L131").Value = 46000
$ws.Range("L24").Value = 45177.5
$ws.Range("L32").Value = 8962.538
$ws.Range("L63").Value = 3500
$ws.Range("L66").Value = 17500
$ws.Range("M122").Value = -8549032
$ws.Range("M32").Value = -8480.237999999999
$ws.Range("M45").Value = -10733.091
$ws.Range("M63").Value = -111113034
$ws.Range("M66").Value = -555565168
$ws.Range("N100").Value = -47341.5
$ws.Range("N122").Value = -10900
$ws.Range("N123").Value = -53657
$ws.Range("N131").Value = -56080
$ws.Range("N24").Value = -45925.5
$ws.Range("N32").Value = -9536.538
$ws.Range("N63").Value = -4872
$ws.Range("N66").Value = -24364

# Sheet 3: BSM
$ws = $wb.Worksheets.Item(3)
$ws.Range("H122").Value = 30499
$ws.Range("J122").Value = 30499
$ws.Range("L122").Value = 30499
$ws.Range("N122").Value = -40299

# Sheet 4: CRP
$ws = $wb.Worksheets.Item(4)
$ws.Range("H113").Value = 0
$ws.Range("H132").Value = 2831.8235
$ws.Range("H16").Value = 0
$ws.Range("H31").Value = 8947.637000000001
$ws.Range("H34").Value = 8947.637000000001
$ws.Range("I113").Value = 0
$ws.Range("I132").Value = 2163.7693
$ws.Range("I16").Value = 0
$ws.Range("I31").Value = 1577
$ws.Range("I34").Value = 1577
$ws.Range("I62").Value = 8801.25
$ws.Range("I65").Value = 8801.25
$ws.Range("J113").Value = 0
$ws.Range("J132").Value = 5003
$ws.Range("J16").Value = 0
$ws.Range("J31").Value = 15884.706
$ws.Range("J34").Value = 15884.706
$ws.Range("J62").Value = 0
$ws.Range("J65").Value = 0
$ws.Range("K113").Value = 0
$ws.Range("K132").Value = 6491.3079
$ws.Range("K16").Value = 0
$ws.Range("K31").Value = 1577
$ws.Range("K34").Value = 1577
$ws.Range("K62").Value = 8801.25
$ws.Range("K65").Value = 44006.25
$ws.Range("L113").Value = 0
$ws.Range("L132").Value = 15009
$ws.Range("L16").Value = 0
$ws.Range("L31").Value = 15884.706
$ws.Range("L34").Value = 15884.706
$ws.Range("L62").Value = 0
$ws.Range("L65").Value = 0
$ws.Range("M113").ClearContents()
$ws.Range("M132").Value = -3961.3079
$ws.Range("M16").ClearContents()
$ws.Range("M31").Value = -1282
$ws.Range("M34").Value = -1375
$ws.Range("M62").Value = -8177.25
$ws.Range("M65").Value = -40886.25
$ws.Range("N113").ClearContents()
$ws.Range("N132").Value = -20069
$ws.Range("N16").ClearContents()
$ws.Range("N31").Value = -16474.706
$ws.Range("N34").Value = -16288.706
$ws.Range("N62").ClearContents()
$ws.Range("N65").ClearContents()

# Sheet 5: CUL
$ws = $wb.Worksheets.Item(5)
$ws.Range("H131").Value = 18966872
$ws.Range("I131").Value = 6250659
$ws.Range("J131").Value = 23811144
$ws.Range("K131").Value = 18751977
$ws.Range("L131").Value = 71433432
$ws.Range("M131").Value = -18746937
$ws.Range("N131").Value = -71443512

# Sheet 6: GSM
$ws = $wb.Worksheets.Item(6)
$ws.Range("H122").Value = 3814947
$ws.Range("H131").Value = 39800
$ws.Range("H69").Value = 42200
$ws.Range("H72").Value = 42200
$ws.Range("I122").Value = 4987623
$ws.Range("J122").Value = 3749.75
$ws.Range("J131").Value = 39800
$ws.Range("J69").Value = 42200
$ws.Range("J72").Value = 42200
$ws.Range("K122").Value = 14962869
$ws.Range("L122").Value = 11249.25
$ws.Range("L131").Value = 39800
$ws.Range("L69").Value = 42200
$ws.Range("L72").Value = 126600
$ws.Range("M122").Value = -14960419
$ws.Range("N122").Value = -16149.25
$ws.Range("N131").Value = -49880
$ws.Range("N69").Value = -43698
$ws.Range("N72").Value = -134088

# Sheet 7: LTW
$ws = $wb.Worksheets.Item(7)
$ws.Range("H122").Value = 10181570
$ws.Range("H132").Value = 10573979
$ws.Range("H16").Value = 547.9091
$ws.Range("I122").Value = 10207509
$ws.Range("I132").Value = 17339266
$ws.Range("J122").Value = 10000000
$ws.Range("J132").Value = 3217.4375
$ws.Range("J16").Value = 0
$ws.Range("K122").Value = 30622527
$ws.Range("K132").Value = 52017798
$ws.Range("L122").Value = 30000000
$ws.Range("L132").Value = 9652.3125
$ws.Range("L16").Value = 0
$ws.Range("M122").Value = -30620077
$ws.Range("M132").Value = -52015268
$ws.Range("N122").Value = -30004900
$ws.Range("N132").Value = -14712.3125
$ws.Range("N16").ClearContents()

# Sheet 8: WVR
$ws = $wb.Worksheets.Item(8)
$ws.Range("H122").Value = 4113
$ws.Range("H123").Value = 29834.285
$ws.Range("I122").Value = 1580
$ws.Range("J122").Value = 8334.666999999999
$ws.Range("J123").Value = 29834.285
$ws.Range("K122").Value = 4740
$ws.Range("L122").Value = 25004.001
$ws.Range("L123").Value = 29834.285
$ws.Range("M122").Value = -2290
$ws.Range("N122").Value = -29904.001
$ws.Range("N123").Value = -39634.285
